$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    if ($Text -match "^-?\d+(\.\d+)?$") {
        # Looks like a plain number -- force text entry the way Excel
        # itself would (leading apostrophe), then strip the resulting
        # quote-prefix style so formatting matches a plain text cell.
        $Cell.Value = "'" + $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

Set-TextCell $ws.Range("D2") "64.028.31"
Set-TextCell $ws.Range("E2") "  -0.71%  "

Set-TextCell $ws.Range("D3") "3.063.26"
Set-TextCell $ws.Range("E3") "  -0.66%  "

Set-TextCell $ws.Range("E4") "  -0.01%  "

Set-TextCell $ws.Range("D5") "561.02"
Set-TextCell $ws.Range("E5") "  +1.22%  "

Set-TextCell $ws.Range("D6") "143.68"
Set-TextCell $ws.Range("E6") "  +0.20%  "

Set-TextCell $ws.Range("E7") "  +0.09%  "

Set-TextCell $ws.Range("D8") "3.061.55"
Set-TextCell $ws.Range("E8") "  -0.58%  "

Set-TextCell $ws.Range("D9") "0.515"
Set-TextCell $ws.Range("E9") "  +4.04%  "

Set-TextCell $ws.Range("E10") "  +1.76%  "

Set-TextCell $ws.Range("D11") "6.15"
Set-TextCell $ws.Range("E11") "  -11.44%  "

Set-TextCell $ws.Range("E12") "  +8.70%  "

Set-TextCell $ws.Range("E13") "  +2.23%  "

Set-TextCell $ws.Range("D14") "35.80"
Set-TextCell $ws.Range("E14") "  +1.44%  "

Set-TextCell $ws.Range("D15") "3.563.10"
Set-TextCell $ws.Range("E15") "  -0.34%  "

Set-TextCell $ws.Range("D16") "64.056.66"
Set-TextCell $ws.Range("E16") "  -0.71%  "

Set-TextCell $ws.Range("D17") "3.060.66"
Set-TextCell $ws.Range("E17") "  -0.65%  "

Set-TextCell $ws.Range("D18") "0.109"
Set-TextCell $ws.Range("E18") "  +1.07%  "

Set-TextCell $ws.Range("D19") "6.82"
Set-TextCell $ws.Range("E19") "  +1.43%  "

Set-TextCell $ws.Range("D20") "478.49"
Set-TextCell $ws.Range("E20") "  -0.40%  "

Set-TextCell $ws.Range("D21") "14.02"
Set-TextCell $ws.Range("E21") "  +2.68%  "

Set-TextCell $ws.Range("D22") "0.686"
Set-TextCell $ws.Range("E22") "  +2.28%  "

Set-TextCell $ws.Range("B23") "InternetComputer(DFINITY)"
Set-TextCell $ws.Range("C23") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D23") "14.42"
Set-TextCell $ws.Range("E23") "  +9.77%  "

Set-TextCell $ws.Range("B24") "Uniswap"
Set-TextCell $ws.Range("C24") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws.Range("D24") "7.58"
Set-TextCell $ws.Range("E24") "  +1.05%  "

Set-TextCell $ws.Range("D25") "82.69"
Set-TextCell $ws.Range("E25") "  +2.52%  "

Set-TextCell $ws.Range("D26") "0.999"
Set-TextCell $ws.Range("E26") "  -0.63%  "

Set-TextCell $ws.Range("D27") "2.82"
Set-TextCell $ws.Range("E27") "  +0.60%  "

Set-TextCell $ws.Range("D28") "8.11"
Set-TextCell $ws.Range("E28") "  +2.57%  "

Set-TextCell $ws.Range("D29") "2.04"
Set-TextCell $ws.Range("E29") "  -0.55%  "

Set-TextCell $ws.Range("D30") "0.999"
Set-TextCell $ws.Range("E30") "  +0.01%  "

Set-TextCell $ws.Range("D31") "26.37"
Set-TextCell $ws.Range("E31") "  +0.87%  "

Set-TextCell $ws.Range("E32") "  -0.76%  "

Set-TextCell $ws.Range("D33") "2.47"
Set-TextCell $ws.Range("E33") "  +1.44%  "

Set-TextCell $ws.Range("E34") "  +1.34%  "

Set-TextCell $ws.Range("D35") "6.26"
Set-TextCell $ws.Range("E35") "  +3.02%  "

Set-TextCell $ws.Range("D36") "54.66"
Set-TextCell $ws.Range("E36") "  -1.06%  "

Set-TextCell $ws.Range("E37") "  +1.40%  "

Set-TextCell $ws.Range("D38") "451.10"
Set-TextCell $ws.Range("E38") "  -3.11%  "

Set-TextCell $ws.Range("E39") "  -1.26%  "

Set-TextCell $ws.Range("E40") "  +5.58%  "

Set-TextCell $ws.Range("D41") "3.029.57"
Set-TextCell $ws.Range("E41") "  +0.43%  "

Set-TextCell $ws.Range("D42") "8.31"
Set-TextCell $ws.Range("E42") "  +0.44%  "

Set-TextCell $ws.Range("E43") "  -1.33%  "

Set-TextCell $ws.Range("E44") "  +3.78%  "

Set-TextCell $ws.Range("D45") "27.95"
Set-TextCell $ws.Range("E45") "  +0.29%  "

Set-TextCell $ws.Range("D46") "2.25"
Set-TextCell $ws.Range("E46") "  +8.80%  "

Set-TextCell $ws.Range("E47") "  +0.00%  "

Set-TextCell $ws.Range("E48") "  +2.00%  "

Set-TextCell $ws.Range("D49") "119.10"

Set-TextCell $ws.Range("D50") "0.0₃0519"
Set-TextCell $ws.Range("E50") "  +0.92%  "

Set-TextCell $ws.Range("D51") "2.11"
Set-TextCell $ws.Range("E51") "  +2.05%  "
